# Insert a new weekly data row at row 4 (pushes existing data rows 4..35 down to 5..36)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new weekly record.
$ws.Cells.Item(4, 1).Value = 7
$ws.Cells.Item(4, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(4, 3).Value = "Ñuble"
$ws.Cells.Item(4, 4).Value = 44761
$ws.Cells.Item(4, 5).Value = 16
$ws.Cells.Item(4, 6).Value = 100112001
$ws.Cells.Item(4, 7).Value = "Berenjena"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 60
$ws.Cells.Item(4, 11).Value = 12000
$ws.Cells.Item(4, 12).Value = 13000
$ws.Cells.Item(4, 13).Value = 12500
$ws.Cells.Item(4, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(4, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(4, 16).Value = 208
$ws.Cells.Item(4, 17).Value = 60
$ws.Cells.Item(4, 18).Value = "Hortaliza"
